$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and 1h Volume-change (E) columns with the latest scraped
# crypto data. A handful of Price values are plain decimals (e.g. "60.04"), which
# Excel would otherwise auto-convert to a Number on assignment; briefly switch those
# specific cells to Text format, write the value, then restore General formatting so
# they keep the same look as every other (already-text) Price cell.

$ws.Range('D2').Value = '37.126.73'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.050.69'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '249.25'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.04'
$ws.Range('D7').NumberFormat = "General"
$ws.Range('E7').Value = '  +8.58%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.389'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0794'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.11'
$ws.Range('D12').NumberFormat = "General"
$ws.Range('E12').Value = '  +7.17%  '
$ws.Range('D13').Value = '2.348.27'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.836'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.82'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('E15').Value = '  +10.09%  '
$ws.Range('D16').Value = '2.046.54'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.34'
$ws.Range('D17').NumberFormat = "General"
$ws.Range('E17').Value = '  +28.53%  '
$ws.Range('D18').Value = '37.137.47'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '76.34'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('E19').Value = '  +3.90%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  -3.92%  '
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '238.78'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('E25').Value = '  +10.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.45'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '169.09'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.25'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.15'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('E30').Value = '  +7.91%  '
$ws.Range('E31').Value = '  +5.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0631'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.65'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('E33').Value = '  +6.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0886'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.24'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.22'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('E40').Value = '  +15.10%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.15'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('E41').Value = '  +18.17%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.69'
$ws.Range('D43').NumberFormat = "General"
$ws.Range('E43').Value = '  -1.78%  '
$ws.Range('E44').Value = '  -0.31%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '97.30'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.51'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.88'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('E47').Value = '  -6.11%  '
$ws.Range('D48').Value = '1.295.96'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.84'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('D51').Value = '2.239.72'
$ws.Range('E51').Value = '  -0.36%  '
